$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.06373566666666666
$ws.Range("H2").Value2 = 0.191207
$ws.Range("I2").Value2 = 0.01058875298517695
$ws.Range("J2").Value2 = 0.01058875298517695
$ws.Range("M2").Value2 = 20.854426
$ws.Range("N2").Value2 = 62.563278
$ws.Range("O2").Value2 = 0.1507164072139519
$ws.Range("P2").Value2 = 0.1507164072139519
$ws.Range("Q2").Value2 = 1.329170744060667
$ws.Range("R2").Value2 = 11.962536696546
$ws.Range("S2").Value2 = 0.001595898806801879
$ws.Range("T2").Value2 = 0.001595898806801879
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.06373566666666666
$ws.Range("H3").Value2 = 0.191207
$ws.Range("I3").Value2 = 0.01058875298517695
$ws.Range("J3").Value2 = 0.01058875298517695
$ws.Range("O3").Value2 = 0.6862909728343718
$ws.Range("P3").Value2 = 0.6862909728343718
$ws.Range("Q3").Value2 = 6.052412606342555
$ws.Range("R3").Value2 = 54.471713457083
$ws.Range("S3").Value2 = 0.00726696558729995
$ws.Range("T3").Value2 = 0.00726696558729995
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.06373566666666666
$ws.Range("H4").Value2 = 0.191207
$ws.Range("I4").Value2 = 0.01058875298517695
$ws.Range("J4").Value2 = 0.01058875298517695
$ws.Range("N4").Value2 = 67.65920700000001
$ws.Range("O4").Value2 = 0.1629926199516763
$ws.Range("P4").Value2 = 0.1629926199516763
$ws.Range("Q4").Value2 = 1.437434888094334
$ws.Range("R4").Value2 = 12.936913992849
$ws.Range("S4").Value2 = 0.001725888591075125
$ws.Range("T4").Value2 = 0.001725888591075125
$ws.Range("I5").Value2 = 0.2961697031425515
$ws.Range("J5").Value2 = 0.2961697031425515
$ws.Range("M5").Value2 = 20.854426
$ws.Range("N5").Value2 = 62.563278
$ws.Range("O5").Value2 = 0.1507164072139519
$ws.Range("P5").Value2 = 0.1507164072139519
$ws.Range("Q5").Value2 = 37.17719218167533
$ws.Range("R5").Value2 = 334.594729635078
$ws.Range("S5").Value2 = 0.04463763358326805
$ws.Range("T5").Value2 = 0.04463763358326805
$ws.Range("I6").Value2 = 0.2961697031425515
$ws.Range("J6").Value2 = 0.2961697031425515
$ws.Range("O6").Value2 = 0.6862909728343718
$ws.Range("P6").Value2 = 0.6862909728343718
$ws.Range("S6").Value2 = 0.2032585936937688
$ws.Range("T6").Value2 = 0.2032585936937688
$ws.Range("I7").Value2 = 0.2961697031425515
$ws.Range("J7").Value2 = 0.2961697031425515
$ws.Range("N7").Value2 = 67.65920700000001
$ws.Range("O7").Value2 = 0.1629926199516763
$ws.Range("P7").Value2 = 0.1629926199516763
$ws.Range("Q7").Value2 = 40.20536362398968
$ws.Range("S7").Value2 = 0.04827347586551468
$ws.Range("T7").Value2 = 0.04827347586551468
$ws.Range("I8").Value2 = 0.6932415438722715
$ws.Range("J8").Value2 = 0.6932415438722715
$ws.Range("M8").Value2 = 20.854426
$ws.Range("N8").Value2 = 62.563278
$ws.Range("O8").Value2 = 0.1507164072139519
$ws.Range("P8").Value2 = 0.1507164072139519
$ws.Range("Q8").Value2 = 87.02029218854932
$ws.Range("R8").Value2 = 783.182629696944
$ws.Range("S8").Value2 = 0.104482874823882
$ws.Range("T8").Value2 = 0.104482874823882
$ws.Range("I9").Value2 = 0.6932415438722715
$ws.Range("J9").Value2 = 0.6932415438722715
$ws.Range("O9").Value2 = 0.6862909728343718
$ws.Range("P9").Value2 = 0.6862909728343718
$ws.Range("S9").Value2 = 0.4757654135533031
$ws.Range("T9").Value2 = 0.4757654135533031
$ws.Range("I10").Value2 = 0.6932415438722715
$ws.Range("J10").Value2 = 0.6932415438722715
$ws.Range("N10").Value2 = 67.65920700000001
$ws.Range("O10").Value2 = 0.1629926199516763
$ws.Range("P10").Value2 = 0.1629926199516763
$ws.Range("Q10").Value2 = 94.10830363437067
$ws.Range("R10").Value2 = 846.9747327093361
$ws.Range("S10").Value2 = 0.1129932554950865
$ws.Range("T10").Value2 = 0.1129932554950865
